# Append the "Team Members:" section (one blank line, an underlined
# heading, and four team-member name paragraphs) after the existing body
# text. We build the new paragraphs as raw WordprocessingML and insert
# them with Range.InsertXML so each paragraph gets exactly the formatting
# we want (Cambria 12pt) instead of inheriting the bullet-list formatting
# of the paragraph currently at the end of the document.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPrPlain = '<w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rPrUnderline = '<w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr>'

$blankPara = "<w:p $wNs><w:pPr>$rPrPlain</w:pPr></w:p>"

$headingPara = "<w:p $wNs><w:pPr>$rPrUnderline</w:pPr><w:r>$rPrUnderline<w:t>Team Members:</w:t></w:r></w:p>"

$names = @("Bunmi Olakanmi", "Chris Inalsingh", "Moné-Renata Holder Seale", "Nihal Joshua")
$namePars = ""
foreach ($name in $names) {
    $namePars += "<w:p $wNs><w:pPr>$rPrPlain</w:pPr><w:r>$rPrPlain<w:t>$name</w:t></w:r></w:p>"
}

$xml = $blankPara + $headingPara + $namePars

$r = $d.Content
$r.Collapse(0)
$r.InsertXML($xml)

Write-Output "Inserted team members section"
